$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("E1").Value = "url_license"
$ws.Range("F1").Value = "text_license"
$ws.Range("G1").Value = "url_source"

# Row 2 (excel_a)
$ws.Range("F2").Value = "This is a license text for entry a"

# Row 3 (excel_b)
$ws.Range("F3").Value = "This is a license text for entry b"

# Row 4 (excel_c)
$ws.Range("F4").Value = "This is a license text for entry c"

# Hyperlinks (also set the displayed text + hyperlink style)
$ws.Hyperlinks.Add($ws.Range("E2"), "https://test.a.html")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://test.b.html")
$ws.Hyperlinks.Add($ws.Range("G3"), "www.source.source.excel_b.xlsx")

# Selection as recorded in the saved file
$ws.Range("E3").Select()
